# Preprocess added to treat instances of NO and NaN
#
# 1) "vomitos" sheet: the categories "NO" and "No" were really the same
#    category (a data-cleaning artifact). Merge the "NO" row's B/C/D counts
#    (the only non-zero ones) into the "No" row, keep the "No" row's E count,
#    drop the now-redundant "NO" row, and shift "Persistente"/"Si" up.
#
# 2) "prueba_torniquete" sheet: rows that previously had a missing/NaN value
#    are now explicitly bucketed into a new "NA" category, inserted before
#    "Negativa".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# vomitos: merge "NO" (row 2) into "No" (row 3), then shift remaining rows up
# ---------------------------------------------------------------------------
$wsVomitos = $wb.Worksheets.Item("vomitos")

# Row 2 becomes "No" with counts merged from old "NO" (B/C/D) and old "No" (E)
$wsVomitos.Cells.Item(2, 1).Value2 = "No"
$wsVomitos.Cells.Item(2, 2).Value2 = 241
$wsVomitos.Cells.Item(2, 3).Value2 = 1697
$wsVomitos.Cells.Item(2, 4).Value2 = 625
$wsVomitos.Cells.Item(2, 5).Value2 = 3670

# Row 3 becomes the old "Persistente" row (was row 4)
$wsVomitos.Cells.Item(3, 1).Value2 = "Persistente"
$wsVomitos.Cells.Item(3, 2).Value2 = 745
$wsVomitos.Cells.Item(3, 3).Value2 = 180
$wsVomitos.Cells.Item(3, 4).Value2 = 382
$wsVomitos.Cells.Item(3, 5).Value2 = 477

# Row 4 becomes the old "Si" row (was row 5)
$wsVomitos.Cells.Item(4, 1).Value2 = "Si"
$wsVomitos.Cells.Item(4, 2).Value2 = 10
$wsVomitos.Cells.Item(4, 3).Value2 = 103
$wsVomitos.Cells.Item(4, 4).Value2 = 35
$wsVomitos.Cells.Item(4, 5).Value2 = 1835

# Drop the now-empty trailing row 5 (shifts rows up / shrinks used range)
$wsVomitos.Rows.Item(5).Delete()

# ---------------------------------------------------------------------------
# prueba_torniquete: insert a new "NA" row before "Negativa"
# ---------------------------------------------------------------------------
$wsTorniquete = $wb.Worksheets.Item("prueba_torniquete")

# Shift existing data rows down by one to make room for the new "NA" row
$wsTorniquete.Cells.Item(4, 1).Value2 = $wsTorniquete.Cells.Item(3, 1).Value2
$wsTorniquete.Cells.Item(4, 2).Value2 = $wsTorniquete.Cells.Item(3, 2).Value2
$wsTorniquete.Cells.Item(4, 3).Value2 = $wsTorniquete.Cells.Item(3, 3).Value2
$wsTorniquete.Cells.Item(4, 4).Value2 = $wsTorniquete.Cells.Item(3, 4).Value2
$wsTorniquete.Cells.Item(4, 5).Value2 = $wsTorniquete.Cells.Item(3, 5).Value2

$wsTorniquete.Cells.Item(3, 1).Value2 = $wsTorniquete.Cells.Item(2, 1).Value2
$wsTorniquete.Cells.Item(3, 2).Value2 = $wsTorniquete.Cells.Item(2, 2).Value2
$wsTorniquete.Cells.Item(3, 3).Value2 = $wsTorniquete.Cells.Item(2, 3).Value2
$wsTorniquete.Cells.Item(3, 4).Value2 = $wsTorniquete.Cells.Item(2, 4).Value2
$wsTorniquete.Cells.Item(3, 5).Value2 = $wsTorniquete.Cells.Item(2, 5).Value2

# New "NA" row goes into row 2
$wsTorniquete.Cells.Item(2, 1).Value2 = "NA"
$wsTorniquete.Cells.Item(2, 2).Value2 = 244
$wsTorniquete.Cells.Item(2, 3).Value2 = 497
$wsTorniquete.Cells.Item(2, 4).Value2 = 262
$wsTorniquete.Cells.Item(2, 5).Value2 = 5067

# Re-apply the header's (bold/bordered/centered) format to column A of all
# data rows, since the row shift above only moved values, not formats.
$wsTorniquete.Range("A1").Copy()
$wsTorniquete.Range("A2:A4").PasteSpecial(-4122)
